$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet view previously had O13 selected; reset the selection to A1
# (the saved file no longer pins the old O13 selection).
$ws.Range("A1").Select() | Out-Null

# Row 13: a new "gave" row (slot order shifted left into E:I) replaces the
# old "tell" row that used E, I, J, K, L.
$ws.Range("E13").Value = "gave"
$ws.Range("F13").Value = "he"
$ws.Range("G13").Value = "gave"
$ws.Range("H13").Value = "me"
$ws.Range("I13").Value = "a message"
$ws.Range("J13").ClearContents() | Out-Null
$ws.Range("K13").ClearContents() | Out-Null
$ws.Range("L13").ClearContents() | Out-Null

# Row 14
$ws.Range("E14").Value = "gave"
$ws.Range("F14").Value = "she"
$ws.Range("G14").Value = "gave"
$ws.Range("H14").Value = "him"
$ws.Range("I14").Value = "a money"
$ws.Range("J14").ClearContents() | Out-Null
$ws.Range("K14").ClearContents() | Out-Null
$ws.Range("L14").ClearContents() | Out-Null

# Row 15
$ws.Range("E15").Value = "gave"
$ws.Range("F15").Value = "Tom"
$ws.Range("G15").Value = "gave"
$ws.Range("H15").Value = "her"
$ws.Range("I15").Value = "ticket"
$ws.Range("J15").ClearContents() | Out-Null
$ws.Range("K15").ClearContents() | Out-Null
$ws.Range("L15").ClearContents() | Out-Null

# Row 16
$ws.Range("E16").Value = "gave"
$ws.Range("F16").Value = "I"
$ws.Range("G16").Value = "gave"
$ws.Range("H16").Value = "Tom"
$ws.Range("I16").Value = "that"
$ws.Range("J16").ClearContents() | Out-Null
$ws.Range("K16").ClearContents() | Out-Null
$ws.Range("L16").ClearContents() | Out-Null
